$wb = $excel.ActiveWorkbook

# Add a brand-new sheet; Excel inserts it *before* the active sheet, so it
# becomes the new first sheet / active tab, pushing the existing
# "marker_template" sheet to position 2.
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "marker_set_template"

# Grab the (now second) original sheet AFTER the insert, since sheet
# references resolve by position.
$oldSheet = $wb.Worksheets.Item(2)

# Populate the new "marker_set_template" sheet: one row describing the
# "Default" marker set.
$newSheet.Range("A1").Value = "id"
$newSheet.Range("B1").Value = "name"
$newSheet.Cells.Item(2, 1).Value = 1
$newSheet.Cells.Item(2, 2).Value = "Default"

# Update the original "marker_template" sheet: column B now stores the
# numeric id of the related marker set instead of the literal set name, and
# its header is renamed accordingly.
$oldSheet.Range("B1").Value = "id_marker_template"
for ($row = 2; $row -le 7; $row++) {
    $oldSheet.Cells.Item($row, 2).Value = 1
}

$oldSheet.Range("E7").Select() | Out-Null
$newSheet.Activate() | Out-Null
$newSheet.Range("E7").Select() | Out-Null
